$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$desc = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

$ws.Range("A10").Value = "JD_009"
$ws.Range("B10").Value = "Senior Test Engineer"
$ws.Range("C10").Value = $desc
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 8

$ws.Rows.Item(10).AutoFit()
